# The commit swaps the two theme parts in the package:
#   ppt/theme/theme1.xml (the deck's main theme, used by the slide master)
#   goes from the custom "Integral" / "Red Violet" color scheme to the
#   stock "Office Theme" / "Office" color scheme (the colors that used to
#   live in ppt/theme/theme2.xml, the notes-master theme).
#
# The font scheme (Arial-based "Office" fonts) and the format scheme
# (fill/line/effect/background styles) are byte-identical between the two
# theme parts already, so the only real content change needed is the 12
# color-scheme entries (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
#
# PowerPoint's object model exposes those via
# Master.Theme.ThemeColorScheme.Colors(index).RGB (VBA/COM RGB() is
# 0xBBGGRR, i.e. the reverse byte order of the "RRGGBB" hex used in the
# OOXML <a:srgbClr val="RRGGBB"/>).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme

# index -> (theme field, target "Office" RRGGBB, COM RGB value)
$scheme.Colors(1).RGB  = 0         # dk1      000000
$scheme.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$scheme.Colors(3).RGB  = 6968388   # dk2      44546A
$scheme.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$scheme.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$scheme.Colors(6).RGB  = 3243501   # accent2  ED7D31
$scheme.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$scheme.Colors(8).RGB  = 49407     # accent4  FFC000
$scheme.Colors(9).RGB  = 12874308  # accent5  4472C4
$scheme.Colors(10).RGB = 4697456   # accent6  70AD47
$scheme.Colors(11).RGB = 12673797  # hlink    0563C1
$scheme.Colors(12).RGB = 7491477   # folHlink 954F72
